$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the material item (header cell) from the cut-data "database"
# and pull its current value into a variable.
$ws.Range("A1").Select()
$material = $ws.Range("A1").Value()

# The sheet only tracks one material, so simplify the label held in the
# variable and write it back to the selected cell.
$material = "STEEL"
$ws.Range("A1").Value = $material
